$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "52.602.52"
$ws.Cells.Item(2, 5).Value = "  +1.78%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.123.46"
$ws.Cells.Item(3, 5).Value = "  +2.08%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.08%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "395.85"
$ws.Cells.Item(5, 5).Value = "  +2.79%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "104.39"
$ws.Cells.Item(6, 5).Value = "  +0.86%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "0.541"
$ws.Cells.Item(7, 5).Value = "  -0.69%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.06%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.604"
$ws.Cells.Item(9, 5).Value = "  +3.08%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "38.17"
$ws.Cells.Item(10, 5).Value = "  +2.77%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +0.96%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "0.0863"
$ws.Cells.Item(12, 5).Value = "  -0.31%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "3.611.80"
$ws.Cells.Item(13, 5).Value = "  +1.90%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "18.86"
$ws.Cells.Item(14, 5).Value = "  +0.78%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "7.88"
$ws.Cells.Item(15, 5).Value = "  +1.39%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  +7.19%  "

# Row 17
$ws.Cells.Item(17, 2).Value = "WrappedEther"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(17, 4).Value = "3.094.91"
$ws.Cells.Item(17, 5).Value = "  +0.96%  "

# Row 18
$ws.Cells.Item(18, 2).Value = "Uniswap"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(18, 4).Value = "11.22"
$ws.Cells.Item(18, 5).Value = "  +6.70%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "52.366.27"
$ws.Cells.Item(19, 5).Value = "  +1.24%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "3.25"
$ws.Cells.Item(20, 5).Value = "  +2.95%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "12.73"
$ws.Cells.Item(21, 5).Value = "  +2.04%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "0.0₃0973"
$ws.Cells.Item(22, 5).Value = "  +0.80%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +0.99%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "269.35"
$ws.Cells.Item(24, 5).Value = "  +0.11%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +1.78%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "8.08"
$ws.Cells.Item(26, 5).Value = "  -4.27%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "27.62"
$ws.Cells.Item(27, 5).Value = "  +2.14%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "7.41"
$ws.Cells.Item(28, 5).Value = "  +1.81%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -2.39%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +0.00%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "0.107"
$ws.Cells.Item(31, 5).Value = "  -0.34%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "10.90"
$ws.Cells.Item(32, 5).Value = "  +5.66%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "36.70"
$ws.Cells.Item(33, 5).Value = "  +6.32%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +9.99%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +0.90%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "49.99"
$ws.Cells.Item(36, 5).Value = "  -0.98%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "0.999"
$ws.Cells.Item(37, 5).Value = "  -0.13%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +0.84%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "NEARProtocol"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(39, 4).Value = "4.10"
$ws.Cells.Item(39, 5).Value = "  +9.77%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "Stacks"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(40, 4).Value = "2.71"
$ws.Cells.Item(40, 5).Value = "  +6.27%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "TheGraph"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(41, 4).Value = "0.294"
$ws.Cells.Item(41, 5).Value = "  +0.54%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "Monero"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(42, 4).Value = "130.65"
$ws.Cells.Item(42, 5).Value = "  +1.19%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "Celestia"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(43, 4).Value = "17.02"
$ws.Cells.Item(43, 5).Value = "  -0.45%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -0.13%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "0.117"
$ws.Cells.Item(45, 5).Value = "  +0.32%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "22.27"
$ws.Cells.Item(46, 5).Value = "  +1.42%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -2.99%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -1.41%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "2.088.92"
$ws.Cells.Item(49, 5).Value = "  +2.04%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "0.0528"
$ws.Cells.Item(50, 5).Value = "  +34.52%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "0.922"
$ws.Cells.Item(51, 5).Value = "  +10.49%  "
